# BookingOutput.xlsx test-result log update (Sreenivas, 10 Oct 2024)
# - Backfills the "Test fail reason" (col D) and "Runtime" (col F) values for a
#   handful of earlier rows that were missing them.
# - Appends newly-run booking test cases as rows 637-660.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Backfill Test fail reason (D) for earlier rows 635-636 ---
$ws.Cells.Item(635,4).Value = "Booking completed"
$ws.Cells.Item(636,4).Value = "Booking completed"

# --- Append new run rows 637-646 ---
# Row 637
$ws.Cells.Item(637,1).Value = "Booking_01"
$ws.Cells.Item(637,2).Value = "ZA12550804"
$ws.Cells.Item(637,3).Value = "Passed"
$ws.Cells.Item(637,4).Value = "Booking completed"
$ws.Cells.Item(637,5).Value = "911907be-622d-4da4-a681-f717782406d6"

# Row 638
$ws.Cells.Item(638,1).Value = "Booking_02"
$ws.Cells.Item(638,2).Value = "ZA12550817"
$ws.Cells.Item(638,3).Value = "Passed"
$ws.Cells.Item(638,4).Value = "Booking completed"
$ws.Cells.Item(638,5).Value = "88d3f0cf-8835-4056-94b7-aab8a8cf9cc6"

# Row 639
$ws.Cells.Item(639,1).Value = "Booking_03"
$ws.Cells.Item(639,2).Value = "ZA12550830"
$ws.Cells.Item(639,3).Value = "Passed"
$ws.Cells.Item(639,4).Value = "Booking completed"
$ws.Cells.Item(639,5).Value = "a075fc87-fe78-4109-a502-a8e86c786ada"

# Row 640
$ws.Cells.Item(640,1).Value = "Booking_04"
$ws.Cells.Item(640,2).Value = "ZA12550838"
$ws.Cells.Item(640,3).Value = "Passed"
$ws.Cells.Item(640,4).Value = "Skipped this test case as this test case is not approved to run"
$ws.Cells.Item(640,5).Value = "76a37c8c-f4cf-46bf-bb6f-6e6c61789288"

# Row 641
$ws.Cells.Item(641,1).Value = "Booking_05"
$ws.Cells.Item(641,2).Value = "NG12550844"
$ws.Cells.Item(641,3).Value = "Passed"
$ws.Cells.Item(641,4).Value = "Skipped this test case as this test case is not approved to run"
$ws.Cells.Item(641,5).Value = "3e3faf3f-595f-4501-ae96-b942f0ecd8d2"

# Row 642
$ws.Cells.Item(642,1).Value = "Booking_01"
$ws.Cells.Item(642,2).Value = "-"
$ws.Cells.Item(642,3).Value = "Skipped"
$ws.Cells.Item(642,4).Value = "Skipped this test case as this test case is not approved to run"
$ws.Cells.Item(642,5).Value = "-"

# Row 643
$ws.Cells.Item(643,1).Value = "Booking_02"
$ws.Cells.Item(643,2).Value = "-"
$ws.Cells.Item(643,3).Value = "Skipped"
$ws.Cells.Item(643,4).Value = "Skipped this test case as this test case is not approved to run"
$ws.Cells.Item(643,5).Value = "-"

# Row 644
$ws.Cells.Item(644,1).Value = "Booking_03"
$ws.Cells.Item(644,2).Value = "-"
$ws.Cells.Item(644,3).Value = "Skipped"
$ws.Cells.Item(644,4).Value = "Skipped this test case as this test case is not approved to run"
$ws.Cells.Item(644,5).Value = "-"

# Row 645
$ws.Cells.Item(645,1).Value = "Booking_04"
$ws.Cells.Item(645,2).Value = "-"
$ws.Cells.Item(645,3).Value = "Skipped"
$ws.Cells.Item(645,4).Value = "Booking completed"
$ws.Cells.Item(645,5).Value = "-"

# Row 646
$ws.Cells.Item(646,1).Value = "Booking_05"
$ws.Cells.Item(646,2).Value = "-"
$ws.Cells.Item(646,3).Value = "Skipped"
$ws.Cells.Item(646,4).Value = "Skipped this test case as this test case is not approved to run"
$ws.Cells.Item(646,5).Value = "-"

# --- Backfill Runtime (F) rows 333-337 ---
$ws.Cells.Item(333,6).Value = "18:51:38"  # Row 333 Runtime backfill
$ws.Cells.Item(334,6).Value = "18:51:39"  # Row 334 Runtime backfill
$ws.Cells.Item(335,6).Value = "18:51:39"  # Row 335 Runtime backfill
$ws.Cells.Item(336,6).Value = "18:51:40"  # Row 336 Runtime backfill
$ws.Cells.Item(337,6).Value = "18:51:40"  # Row 337 Runtime backfill

# --- Append new run row 647 ---
# Row 647
$ws.Cells.Item(647,1).Value = "Booking_06"
$ws.Cells.Item(647,2).Value = "NG12550889"
$ws.Cells.Item(647,3).Value = "Passed"
$ws.Cells.Item(647,4).Value = "Skipped this test case as this test case is not approved to run"
$ws.Cells.Item(647,5).Value = "28cc9b50-b41e-427d-9df3-5b3e23008243"

# --- Backfill Runtime (F) rows 338-342 ---
$ws.Cells.Item(338,6).Value = "18:24:49"  # Row 338 Runtime backfill
$ws.Cells.Item(339,6).Value = "18:24:51"  # Row 339 Runtime backfill
$ws.Cells.Item(340,6).Value = "18:24:52"  # Row 340 Runtime backfill
$ws.Cells.Item(341,6).Value = "18:24:52"  # Row 341 Runtime backfill
$ws.Cells.Item(342,6).Value = "18:24:53"  # Row 342 Runtime backfill

# Row 343: explicit empty cell in column F (style only, no value)
$ws.Cells.Item(343,6).Style = "Normal"

# --- Append new run rows 648-660 ---
# Row 648
$ws.Cells.Item(648,1).Value = "Booking_01"
$ws.Cells.Item(648,2).Value = "-"
$ws.Cells.Item(648,3).Value = "Skipped"
$ws.Cells.Item(648,4).Value = "Skipped this test case as this test case is not approved to run"
$ws.Cells.Item(648,5).Value = "-"

# Row 649
$ws.Cells.Item(649,1).Value = "Booking_02"
$ws.Cells.Item(649,2).Value = "-"
$ws.Cells.Item(649,3).Value = "Skipped"
$ws.Cells.Item(649,4).Value = "Skipped this test case as this test case is not approved to run"
$ws.Cells.Item(649,5).Value = "-"

# Row 650
$ws.Cells.Item(650,1).Value = "Booking_03"
$ws.Cells.Item(650,2).Value = "-"
$ws.Cells.Item(650,3).Value = "Skipped"
$ws.Cells.Item(650,4).Value = "Skipped this test case as this test case is not approved to run"
$ws.Cells.Item(650,5).Value = "-"

# Row 651
$ws.Cells.Item(651,1).Value = "Booking_04"
$ws.Cells.Item(651,2).Value = "-"
$ws.Cells.Item(651,3).Value = "Skipped"
$ws.Cells.Item(651,4).Value = "Booking completed"
$ws.Cells.Item(651,5).Value = "-"

# Row 652
$ws.Cells.Item(652,1).Value = "Booking_05"
$ws.Cells.Item(652,2).Value = "-"
$ws.Cells.Item(652,3).Value = "Skipped"
$ws.Cells.Item(652,4).Value = "Booking completed"
$ws.Cells.Item(652,5).Value = "-"

# Row 653
$ws.Cells.Item(653,1).Value = "Booking_01"
$ws.Cells.Item(653,2).Value = "ZA12557861"
$ws.Cells.Item(653,3).Value = "Passed"
$ws.Cells.Item(653,4).Value = "Booking completed"
$ws.Cells.Item(653,5).Value = "dfe1ff76-27ad-4294-9019-7d191438f33e"

# Row 654
$ws.Cells.Item(654,1).Value = "Booking_02"
$ws.Cells.Item(654,2).Value = "ZA12557877"
$ws.Cells.Item(654,3).Value = "Passed"
$ws.Cells.Item(654,4).Value = "Booking completed"
$ws.Cells.Item(654,5).Value = "42d2abb1-ce91-4f1a-86fa-7c348d7b84fa"

# Row 655
$ws.Cells.Item(655,1).Value = "Booking_03"
$ws.Cells.Item(655,2).Value = "ZA12557885"
$ws.Cells.Item(655,3).Value = "Passed"
$ws.Cells.Item(655,4).Value = "Booking completed"
$ws.Cells.Item(655,5).Value = "489ccab0-b335-4f94-afe0-40b5dbe81a30"

# Row 656
$ws.Cells.Item(656,1).Value = "Booking_05"
$ws.Cells.Item(656,2).Value = "NG12557910"
$ws.Cells.Item(656,3).Value = "Passed"
$ws.Cells.Item(656,4).Value = "Booking completed"
$ws.Cells.Item(656,5).Value = "6c47f2df-1089-43ad-b0eb-c2b27511eb41"

# Row 657
$ws.Cells.Item(657,1).Value = "Booking_06"
$ws.Cells.Item(657,2).Value = "NG12557916"
$ws.Cells.Item(657,3).Value = "Passed"
$ws.Cells.Item(657,4).Value = "Booking completed"
$ws.Cells.Item(657,5).Value = "e09b1d56-3307-4bf5-a12e-0973f1a62074"

# Row 658
$ws.Cells.Item(658,1).Value = "Booking_08"
$ws.Cells.Item(658,2).Value = "NG12557921"
$ws.Cells.Item(658,3).Value = "Passed"
$ws.Cells.Item(658,4).Value = "Booking completed"
$ws.Cells.Item(658,5).Value = "60193234-c143-4885-a642-ce17af7377a6"

# Row 659
$ws.Cells.Item(659,1).Value = "Booking_01"
$ws.Cells.Item(659,2).Value = "ZA00182878"
$ws.Cells.Item(659,3).Value = "Passed"
$ws.Cells.Item(659,5).Value = "12db59c9-ef07-4774-adfa-3f59ea47d040"

# Row 660
$ws.Cells.Item(660,1).Value = "Booking_02"
$ws.Cells.Item(660,2).Value = "ZA00182879"
$ws.Cells.Item(660,3).Value = "Passed"
$ws.Cells.Item(660,5).Value = "54ce1b25-4390-427a-b9c9-354a7cff377e"
